$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two extra contract rows (row 3 "hd2" and row 4 "hd3") are removed -
# only the first contract row remains, now updated with the renewed
# ("gia han") contract details.
$ws.Rows("3:4").Delete()

# Plain (non date-like) text values can be assigned directly.
$ws.Range("B2").Value = "HD1"
$ws.Range("C2").Value = "NV005"
$ws.Range("D2").Value = "SV005"
$ws.Range("E2").Value = "A101"
$ws.Range("H2").Value = "Hết hạn"

# F2/G2 hold date-shaped text ("2024-06-24" / "2024-06-26") that must stay
# plain text (same style as the rest of the row), not get auto-converted
# into a real date serial number. Routing the literal through a text
# formula and then flattening it back to a static value (copy/paste as
# values onto itself) keeps it typed as text without touching the cell's
# number format/style.
$ws.Range("F2").Formula = "=""2024-06-24"""
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

$ws.Range("G2").Formula = "=""2024-06-26"""
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)
